$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert new rows (bottom-up so earlier row numbers stay stable) ---
# Original row 19 ("Visitor" / "View site") gets one blank row inserted above it
# (ends up at row 22 once the later 2-row insert above it is applied too).
$ws.Rows("19:19").Insert()
# Two new rows are inserted before original row 16 ("send requests") for the
# new "Add requirements" / "Search users needs" feature rows.
$ws.Rows("16:17").Insert()

# --- New "Features" column (G) filled in row order, matching the rows that
# already existed before this edit ---
$ws.Range("G4").Value = "Features"
$ws.Range("G7").Value = "Admin can view the details of all the registred sellers"
$ws.Range("G8").Value = "admin can view the details of all the registred buyers"
$ws.Range("G9").Value = "view all the property details that can be registred by sellers"
$ws.Range("G10").Value = "view feedbacks from buyers ,sellers and viewers"
$ws.Range("G13").Value = "After registration a fixed amount of fees can be pay by both the sellers and buyers"
$ws.Range("G14").Value = "After paying amount they can login to their personal profile"
$ws.Range("G15").Value = "sellers and buyers can edit ,delete,add details to there profile"
$ws.Range("G18").Value = "When matching requirements come the buyer can send request to the seller"
$ws.Range("G19").Value = "When a  new property is added to the site it can be view by the buyer as a notification"
$ws.Range("G22").Value = "visitor can just view the site and wants to know more details they can create an account"
$ws.Range("G23").Value = "they can send feed backs to admin"

# --- New row 20: send feedbacks ---
$ws.Range("C20").Value = "send feedbacks"
$ws.Range("G20").Value = "both buyer and seller can send feedbacks to admin about the site"

# --- New rows 16 & 17: Add requirements / Search users needs ---
$ws.Range("C16").Value = "Add requirements"
$ws.Range("G16").Value = "They can add there requirements"
$ws.Range("C17").Value = "Search users needs"
$ws.Range("G17").Value = "They are able to search property according to there needs"

# --- Selection matches the final saved state ---
$ws.Range("G17").Select()
